$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1013748.3
$ws.Range("I31").Value = 1013748.3
$ws.Range("K31").Value = 3041244.9
$ws.Range("M31").Value = -3041014.9
$ws.Range("H70").Value = 16767928
$ws.Range("I70").Value = 41917492
$ws.Range("J70").Value = 1550.5
$ws.Range("K70").Value = 125752476
$ws.Range("L70").Value = 4651.5
$ws.Range("M70").Value = -125752206
$ws.Range("N70").Value = -5191.5
$ws.Range("H73").Value = 16767928
$ws.Range("I73").Value = 41917492
$ws.Range("J73").Value = 1550.5
$ws.Range("K73").Value = 125752476
$ws.Range("L73").Value = 4651.5
$ws.Range("M73").Value = -125751540
$ws.Range("N73").Value = -6523.5
$ws.Range("H76").Value = 5133.9287
$ws.Range("I76").Value = 6571
$ws.Range("J76").Value = 3696.8572
$ws.Range("K76").Value = 6571
$ws.Range("L76").Value = 3696.8572
$ws.Range("M76").Value = -6256
$ws.Range("N76").Value = -4326.8572
$ws.Range("H79").Value = 5133.9287
$ws.Range("I79").Value = 6571
$ws.Range("J79").Value = 3696.8572
$ws.Range("K79").Value = 6571
$ws.Range("L79").Value = 3696.8572
$ws.Range("M79").Value = -5479
$ws.Range("N79").Value = -5880.8572
$ws.Range("H92").Value = 2438.5
$ws.Range("I92").Value = 2918
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2918
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -1670
$ws.Range("N92").Value = -3496
$ws.Range("H107").Value = 389.91666
$ws.Range("I107").Value = 287.69232
$ws.Range("J107").Value = 510.72726
$ws.Range("K107").Value = 287.69232
$ws.Range("L107").Value = 510.72726
$ws.Range("M107").Value = 1632.30768
$ws.Range("N107").Value = -4350.72726
$ws.Range("H116").Value = 50001624
$ws.Range("I116").Value = 50001624
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 50001624
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -49998182
$ws.Range("H132").Value = 1622.4929
$ws.Range("I132").Value = 1292.0189
$ws.Range("J132").Value = 2595.5557
$ws.Range("K132").Value = 3876.0567
$ws.Range("L132").Value = 7786.6671
$ws.Range("M132").Value = -1346.0567
$ws.Range("N132").Value = -12846.6671
$ws.Range("H135").Value = 773.087
$ws.Range("I135").Value = 732.4286
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 6591.8574
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -4056.8574
$ws.Range("N135").Value = -15870
$ws.Range("H137").Value = 1347.3334
$ws.Range("I137").Value = 1347.3334
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4042.0002
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -1492.0002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3253.25
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3253.25
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = 3253.25
$ws.Range("N2").Value = -3479.25
$ws.Range("H32").Value = 12620.45
$ws.Range("I32").Value = 14084.98
$ws.Range("J32").Value = 4321.4443
$ws.Range("K32").Value = 14084.98
$ws.Range("L32").Value = 4321.4443
$ws.Range("M32").Value = -13797.98
$ws.Range("N32").Value = -4895.4443
$ws.Range("H45").Value = 1459.4667
$ws.Range("I45").Value = 1124.5
$ws.Range("J45").Value = 1842.2858
$ws.Range("K45").Value = 1124.5
$ws.Range("L45").Value = 1842.2858
$ws.Range("M45").Value = -747.5
$ws.Range("N45").Value = -2596.2858
$ws.Range("H109").Value = 31437.5
$ws.Range("J109").Value = 31437.5
$ws.Range("L109").Value = 31437.5
$ws.Range("N109").Value = -34211.5
$ws.Range("H116").Value = 3253.25
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3253.25
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").Value = 3253.25
$ws.Range("N116").Value = -7841.25
$ws.Range("H122").Value = 5531
$ws.Range("I122").Value = 6062
$ws.Range("K122").Value = 18186
$ws.Range("M122").Value = -15736
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3253.25
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3253.25
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 3253.25
$ws.Range("N3").Value = -3481.25
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("N34").Value = 0
$ws.Range("H107").Value = 59333.223
$ws.Range("I107").Value = 74857
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 74857
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -72937
$ws.Range("N107").Value = -8840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1749.1136
$ws.Range("I31").Value = 1298.4474
$ws.Range("J31").Value = 4603.3335
$ws.Range("K31").Value = 1298.4474
$ws.Range("L31").Value = 4603.3335
$ws.Range("M31").Value = -1003.4474
$ws.Range("N31").Value = -5193.3335
$ws.Range("H34").Value = 1749.1136
$ws.Range("I34").Value = 1298.4474
$ws.Range("J34").Value = 4603.3335
$ws.Range("K34").Value = 1298.4474
$ws.Range("L34").Value = 4603.3335
$ws.Range("M34").Value = -1096.4474
$ws.Range("N34").Value = -5007.3335
$ws.Range("H86").Value = 2817.6667
$ws.Range("I86").Value = 1762.8
$ws.Range("K86").Value = 1762.8
$ws.Range("M86").Value = -639.8
$ws.Range("H89").Value = 2817.6667
$ws.Range("I89").Value = 1762.8
$ws.Range("K89").Value = 8814
$ws.Range("M89").Value = -3198
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0
$ws.Range("H122").Value = 3401.889
$ws.Range("I122").Value = 2297.8
$ws.Range("K122").Value = 6893.400000000001
$ws.Range("M122").Value = -4443.400000000001
$ws.Range("H131").Value = 24325
$ws.Range("J131").Value = 24325
$ws.Range("L131").Value = 24325
$ws.Range("N131").Value = -34405
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1073587.1
$ws.Range("J12").Value = 1136739.2
$ws.Range("L12").Value = 3410217.6
$ws.Range("N12").Value = -3410563.6
$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -2189
$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -4944
$ws.Range("H92").Value = 724.2857
$ws.Range("I92").Value = 693.3333
$ws.Range("K92").Value = 2079.9999
$ws.Range("M92").Value = -831.9998999999998
$ws.Range("H97").Value = 850
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 2700
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -3692
$ws.Range("H107").Value = 1566.6666
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 6000
$ws.Range("M107").Value = -4080
$ws.Range("H129").Value = 3334099.8
$ws.Range("I129").Value = 192.5
$ws.Range("J129").Value = 4546429.5
$ws.Range("K129").Value = 577.5
$ws.Range("L129").Value = 13639288.5
$ws.Range("M129").Value = 4422.5
$ws.Range("N129").Value = -13649288.5
$ws.Range("H131").Value = 13637.167
$ws.Range("I131").Value = 348.7647
$ws.Range("J131").Value = 17744.49
$ws.Range("K131").Value = 1046.2941
$ws.Range("L131").Value = 53233.47
$ws.Range("M131").Value = 3993.7059
$ws.Range("N131").Value = -63313.47
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8053.3687
$ws.Range("I70").Value = 9601.362999999999
$ws.Range("J70").Value = 5924.875
$ws.Range("K70").Value = 9601.362999999999
$ws.Range("L70").Value = 5924.875
$ws.Range("M70").Value = -9331.362999999999
$ws.Range("N70").Value = -6464.875
$ws.Range("H73").Value = 8053.3687
$ws.Range("I73").Value = 9601.362999999999
$ws.Range("J73").Value = 5924.875
$ws.Range("K73").Value = 9601.362999999999
$ws.Range("L73").Value = 5924.875
$ws.Range("M73").Value = -8665.362999999999
$ws.Range("N73").Value = -7796.875
$ws.Range("H102").Value = 2179.4075
$ws.Range("I102").Value = 2478.1538
$ws.Range("J102").Value = 1902
$ws.Range("K102").Value = 2478.1538
$ws.Range("L102").Value = 1902
$ws.Range("M102").Value = -856.1538
$ws.Range("N102").Value = -5146
$ws.Range("H122").Value = 3239.5
$ws.Range("I122").Value = 3341
$ws.Range("J122").Value = 3002.6667
$ws.Range("K122").Value = 10023
$ws.Range("L122").Value = 9008.000100000001
$ws.Range("M122").Value = -7573
$ws.Range("N122").Value = -13908.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1157.05
$ws.Range("I93").Value = 696.3125
$ws.Range("K93").Value = 696.3125
$ws.Range("M93").Value = 551.6875
$ws.Range("H122").Value = 66668036
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 100001060
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 300003180
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -300008080
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 78066.664
$ws.Range("H90").Value = 78066.664
$ws.Range("H100").Value = 8905.68
$ws.Range("I100").Value = 25521.5
$ws.Range("J100").Value = 1086.4706
$ws.Range("K100").Value = 51043
$ws.Range("L100").Value = 2172.9412
$ws.Range("M100").Value = -50502
$ws.Range("N100").Value = -3254.9412
$ws.Range("H123").Value = 39571.855
$ws.Range("J123").Value = 39571.855
$ws.Range("L123").Value = 39571.855
$ws.Range("N123").Value = -49371.855
$ws.Range("H132").Value = 3148.7407
$ws.Range("I132").Value = 2805.25
$ws.Range("J132").Value = 3423.5334
$ws.Range("K132").Value = 8415.75
$ws.Range("L132").Value = 10270.6002
$ws.Range("M132").Value = -5885.75
$ws.Range("N132").Value = -15330.6002
